$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L1 value "Lunch" -> "Fomal"
$ws.Range("L1").Value = "Fomal"

# Update O1 value "Music" -> "Non-Music"
$ws.Range("O1").Value = "Non-Music"

# Set column widths (E and O) to match the target layout
$ws.Columns.Item(5).ColumnWidth = 13.498697916666666
$ws.Columns.Item(15).ColumnWidth = 10.276041666666666

# Update selection to O1
[void]$ws.Range("O1").Select()
